# Auto-generated Word COM-interop script.
# Applies the OOXML diff by replacing whole paragraphs via Range.InsertXML
# with precise run/bookmark structure, located robustly via Find.
$d = $word.ActiveDocument

# 1) Split the trailing run of the 'Contexto:' paragraph and drop the
#    trailing space after 'opciones:'.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("desarrollo las siguientes opciones", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "anchor paragraph 1 not found" }
$rng1.Expand(4) | Out-Null
$x1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00AF3334" w:rsidRPr="00A2772F" w:rsidRDefault="00AF3334" w:rsidP="00AF3334"><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="240" w:afterAutospacing="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="24292E"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="00A2772F"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="24292E"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Contexto: Como se cuenta con una gran cantidad de información acerca de los clientes de la compañía, se realizará un análisis riguroso de sus comportamientos a lo largo del tiempo. Dicho lo anterior, la información de compras y ventas de cada uno serían de gran importancia en el análisis planteado, ya que así se obtendrían las preferencias que tiene el cliente; ya sea de los productos que más frecuenta o los más consultados por clientes similares a él. Esto con el objetivo de sugerir los productos y lograr una conformidad mayor tanto como para el cliente como para la empresa. Por estas razones se sugiere la implementa</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="24292E"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>ción de un software el cual esté</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="24292E"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> en capacidad de proporcionar una predicción lo más acertada posible, teniendo como alternativas de desa</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="24292E"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>rrollo las siguientes opciones:</w:t></w:r></w:p>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$rng1.InsertXML($x1)

# 2) Remove the proofErr spell-check wrapper around 'BruteForce' and move
#    the _GoBack bookmark in between 'Brut' and 'eForce'.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("BruteForce", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "anchor paragraph 2 not found" }
$rng2.Expand(4) | Out-Null
$x2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00AF3334" w:rsidRPr="00A2772F" w:rsidRDefault="00AF3334" w:rsidP="00AF3334"><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="240" w:afterAutospacing="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="24292E"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="24292E"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Alternativa 1: Con el fin de realizar una búsqueda de reglas de asociación, se propone implementar la estrategia “Brut</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="24292E"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">eForce” sobre los datos que se tienen en poder. Este algoritmo consiste en enumerar sistemáticamente todos los posibles candidatos, con el fin de revisar si dicho candidato satisface la solución del problema. </w:t></w:r></w:p>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$rng2.InsertXML($x2)

# 3) Remove the old _GoBack bookmark left on the trailing empty Heading 1
#    paragraph (the bookmark now lives inside 'BruteForce' instead).
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.ParagraphStyle.NameLocal -eq "Heading 1" -and $p.Range.Text.Trim().Length -eq 0) {
        $rng3 = $p.Range
$x3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="005F49AA" w:rsidRPr="007F5006" w:rsidRDefault="005F49AA" w:rsidP="00AF3334"><w:pPr><w:pStyle w:val="Ttulo1"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="24292E"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

        $rng3.InsertXML($x3)
        break
    }
}
